$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K quarterly data to F:M,
# making room for the two most-recently reported quarters).
$ws.Columns("D:E").Insert()

# Copy cell formatting (number/date format, font, alignment) from the data that was
# just shifted into F:G into the newly inserted, still-blank D:E columns so every row
# keeps its correct style (date header rows vs numeric data rows).
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new D (quarter ending 2018-12-31) and E (quarter ending 2018-09-30)
# columns with the latest reported financial figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 4801100
$ws.Range("E8").Value = 4658000
$ws.Range("D9").Value = 4624900
$ws.Range("E9").Value = 4432200
$ws.Range("D10").Value = 176200
$ws.Range("E10").Value = 225800
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 4640400
$ws.Range("E17").Value = 4511400
$ws.Range("D18").Value = 160700
$ws.Range("E18").Value = 146600
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 212400
$ws.Range("E21").Value = 199500
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 160700
$ws.Range("E23").Value = 146500
$ws.Range("D24").Value = 82800
$ws.Range("E24").Value = 50500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 77800
$ws.Range("E26").Value = 96000
$ws.Range("D27").Value = 50200
$ws.Range("E27").Value = 77300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 50200
$ws.Range("E33").Value = 77300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 50200
$ws.Range("E35").Value = 77300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 1764700
$ws.Range("E41").Value = 1679800
$ws.Range("D42").Value = 214800
$ws.Range("E42").Value = 243000
$ws.Range("D43").Value = 1534300
$ws.Range("E43").Value = 1724100
$ws.Range("D44").Value = 1545000
$ws.Range("E44").Value = 1464400
$ws.Range("D45").Value = 382000
$ws.Range("E45").Value = 472900
$ws.Range("D46").Value = 5440900
$ws.Range("E46").Value = 5584200
$ws.Range("D47").Value = 938500
$ws.Range("E47").Value = 859200
$ws.Range("D48").Value = 1013700
$ws.Range("E48").Value = 1029200
$ws.Range("D49").Value = 533600
$ws.Range("E49").Value = 553300
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 986900
$ws.Range("E52").Value = 1122100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 8913600
$ws.Range("E54").Value = 9148000
$ws.Range("D57").Value = 1638900
$ws.Range("E57").Value = 1619200
$ws.Range("D58").Value = 26900
$ws.Range("E58").Value = 56900
$ws.Range("D59").Value = 1886700
$ws.Range("E59").Value = 1985900
$ws.Range("D60").Value = 3552500
$ws.Range("E60").Value = 3662100
$ws.Range("D61").Value = 1661600
$ws.Range("E61").Value = 1667400
$ws.Range("D62").Value = 581500
$ws.Range("E62").Value = 618600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 5950500
$ws.Range("E66").Value = 6105000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 3422200
$ws.Range("E72").Value = 3401500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 2963200
$ws.Range("E76").Value = 3043000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 50200
$ws.Range("E81").Value = 77300
$ws.Range("D83").Value = 51700
$ws.Range("E83").Value = 53000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 173400
$ws.Range("E89").Value = 121300
$ws.Range("D91").Value = -62300
$ws.Range("E91").Value = -38200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 73100
$ws.Range("E94").Value = -155700
$ws.Range("D96").Value = -29500
$ws.Range("E96").Value = -29600
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -132100
$ws.Range("E100").Value = 32100
$ws.Range("D101").Value = -29500
$ws.Range("E101").Value = 200
$ws.Range("D102").Value = 85000
$ws.Range("E102").Value = -2200
